$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8625954198473282
$ws.Range("B3").Value = 0.9244935543278084
$ws.Range("B4").Value = 0.8297520661157025
$ws.Range("B5").Value = 0.8745644599303136
